$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Отчёт о движении"

# Update the reporting period text in row 2 (merged A2:G2) to the new range.
$ws.Range("A2:G2").Value = "Период: 2023-10-01 - 2023-10-31"

# Give the header/data rows (6-10) their own distinct (but visually identical)
# plain bordered style, separate from the blank separator row above (row 5).
$ws.Range("A6:G10").WrapText = $false

# Select the title row, as in the final workbook.
$ws.Range("A1:G1").Select()

$wb.Save()
